$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift Timestamp (column A) by +2 days for every data row (2..193),
#    and rebuild the Lookup (column E) text from the new date + existing Quarter (column D).
for ($r = 2; $r -le 193; $r++) {
    $serial = $ws.Cells.Item($r, 1).Value2
    $newSerial = $serial + 2
    $ws.Cells.Item($r, 1).Value2 = $newSerial
    $quarter = $ws.Cells.Item($r, 4).Value2
    $dt = [DateTime]::FromOADate($newSerial)
    $ws.Cells.Item($r, 5).Value = ($dt.ToString("dd.MM.yyyy") + [string]$quarter)
}

# 2) Apply the updated Notified Production (column B) readings.
$bValues = @{
    22 = 0.466
    23 = 0.474
    24 = 0.482
    25 = 0.5
    26 = 0.862
    27 = 0.98
    28 = 1.426
    29 = 3.429
    30 = 12.362
    31 = 24.669
    32 = 39.247
    33 = 56.694
    34 = 106.549
    35 = 132.406
    36 = 165.312
    37 = 191.688
    38 = 240.953
    39 = 269.027
    40 = 304.853
    41 = 323.894
    42 = 368.259
    43 = 390.89
    44 = 422.862
    45 = 441.691
    46 = 477.804
    47 = 488.208
    48 = 497.381
    49 = 502.15
    50 = 505.301
    51 = 500.954
    52 = 493.74
    53 = 477.834
    54 = 445.929
    55 = 421.228
    56 = 393.544
    57 = 358.952
    58 = 287.644
    59 = 249.545
    60 = 209.071
    61 = 174.115
    62 = 118.123
    63 = 86.131
    64 = 57.959
    65 = 38.028
    66 = 24.885
    67 = 11.434
    68 = 9.19
    69 = 7.096
    70 = 0.659
    73 = 0
    74 = 0.65
    75 = 0
    118 = 0.472
    119 = 0.481
    120 = 0.499
    121 = 0.509
    122 = 1.341
    123 = 1.554
    124 = 1.886
    125 = 2.86
    126 = 19.708
    127 = 31.981
    128 = 47.344
    129 = 65.388
    130 = 105.727
    131 = 121.828
    132 = 152.675
    133 = 181.732
    134 = 232.458
    135 = 261.499
    136 = 291.215
    137 = 318.427
    138 = 372.687
    139 = 396.213
    140 = 422.951
    141 = 441.843
    142 = 475.607
    143 = 486.664
    144 = 493.508
    145 = 495.279
    146 = 507.976
    147 = 501.512
    148 = 492.846
    149 = 477.147
    150 = 445.878
    151 = 421.425
    152 = 392.221
    153 = 357.985
    154 = 286.131
    155 = 248.384
    156 = 202.187
    157 = 165.317
    158 = 105.128
    159 = 78.164
    160 = 53.878
    161 = 33.918
    162 = 21.523
    163 = 8.908
    164 = 7.349
    165 = 6.832
    166 = 5.151
    170 = 4.85
    171 = 3.15
    172 = 2.45
    173 = 0.65
    182 = 0.55
}
foreach ($r in $bValues.Keys) {
    $ws.Cells.Item([int]$r, 2).Value2 = $bValues[$r]
}

# 3) Apply the updated Actual Production (column C) readings.
$cValues = @{
    30 = 0
    31 = 13
    32 = 31
    33 = 57
    34 = 91
    35 = 127
    36 = 174
    37 = 217
    38 = 260
    39 = 297
    40 = 329
    41 = 371
    42 = 418
    43 = 407
    44 = 428
    45 = 455
    46 = 496
    47 = 490
    48 = 532
    49 = 527
    50 = 514
    51 = 496
    52 = 492
    53 = 464
    54 = 452
    55 = 443
    56 = 412
    57 = 372
    58 = 351
    59 = 295
    60 = 234
    61 = 171
    62 = 125
    63 = 87
    64 = 58
    65 = 19
    66 = 3
}
foreach ($r in $cValues.Keys) {
    $ws.Cells.Item([int]$r, 3).Value2 = $cValues[$r]
}
